# JS-Frameworks-Self-Evaluation-Protocol.xlsx — "Fixed Self Evaluation Protocol"
#
# Changes:
#  - GitHub commit-days score (C8): 8 -> 100
#  - GitHub commit-count score (C9): 25 -> (cleared/blank)
#  - Total score formula in C32 recalculates accordingly (254 -> 321)
#  - The sheet's scroll position / selection moved (view was left on
#    E11 with the window scrolled down to row 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- data edits -----------------------------------------------------
$ws.Range("C8").Value = 100
$ws.Range("C9").ClearContents()

# --- view / selection state ------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$excel.Goto($ws.Range("E11"), $false)

$wb.Save()
